# Updates the "cryptos" price/volume table to the latest scraped values
# (GitHub Actions symbol-list refresh). D/E columns hold numeric-looking
# text (prices / percentages) that must stay stored as text, so each is
# written with a leading apostrophe (forces text entry, same as typing
# it into Excel) and the cell style is reset back to "Normal" right
# after so no stray number-format/quote-prefix style sticks around.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'326.56"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'-1.27%"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'44.12"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'-0.90%"
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.Value = "'5.291"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'-4.42%"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'0.08334"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'2.17%"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'1.939"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'-5.68%"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'-0.69%"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'-3.93%"
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = "'0.1124"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'0.38%"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'0.1891"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'-0.37%"
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = "'0.09660"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'-3.58%"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'0.04597"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'-2.92%"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'0.16%"
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'0.001292"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'2.52%"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'0.005872"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'-2.32%"
$c.Style = "Normal"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c = $ws.Range("D16")
$c.Value = "'3.402"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'1.64%"
$c.Style = "Normal"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$c = $ws.Range("D17")
$c.Value = "'4.406"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'-0.44%"
$c.Style = "Normal"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$c = $ws.Range("D18")
$c.Value = "'0.3357"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'0.19%"
$c.Style = "Normal"
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$c = $ws.Range("D19")
$c.Value = "'8.543"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'-16.53%"
$c.Style = "Normal"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$c = $ws.Range("D20")
$c.Value = "'0.1371"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'-1.31%"
$c.Style = "Normal"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$c = $ws.Range("D21")
$c.Value = "'0.2577"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'0.44%"
$c.Style = "Normal"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$c = $ws.Range("D22")
$c.Value = "'0.04158"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'1.45%"
$c.Style = "Normal"
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$c = $ws.Range("D23")
$c.Value = "'0.001234"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'-5.25%"
$c.Style = "Normal"
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$c = $ws.Range("D24")
$c.Value = "'0.004410"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'0.24%"
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.Value = "'0.0002984"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'-19.99%"
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.Value = "'0.02675"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'-0.44%"
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'0.05565"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'-1.96%"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'0.007847"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'2.98%"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'-0.87%"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'0.007325"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'-2.38%"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'0.007841"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'-5.38%"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'0.3505"
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.Value = "'0.00006860"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'-2.60%"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'0.44%"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'0.003498"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'-0.95%"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'0.003536"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'40.70%"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'0.00002104"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'0.44%"
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.Value = "'0.0002004"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'0.44%"
$c.Style = "Normal"
